$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (bold/bordered/centered style used by A1:H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-48 for columns I (I0) and J (IF)
$data = @(
    @(2, 7, 8),
    @(3, 8, 8),
    @(4, 8, 8),
    @(5, 7, 7),
    @(6, 7, 7),
    @(7, 9, 9),
    @(8, 8, 8),
    @(9, 6, 7),
    @(10, 6, 6),
    @(11, 8, 8),
    @(12, 9, 9),
    @(13, 8, 8),
    @(14, 8, 8),
    @(15, 7, 7),
    @(16, 9, 9),
    @(17, 9, 9),
    @(18, 8, 8),
    @(19, 8, 8),
    @(20, 7, 7),
    @(21, 8, 8),
    @(22, 8, 8),
    @(23, 9, 9),
    @(24, 9, 9),
    @(25, 7, 7),
    @(26, 8, 9),
    @(27, 7, 7),
    @(28, 9, 9),
    @(29, 8, 8),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 8, 8),
    @(33, 9, 9),
    @(34, 8, 8),
    @(35, 9, 9),
    @(36, 7, 7),
    @(37, 8, 8),
    @(38, 9, 9),
    @(39, 9, 9),
    @(40, 8, 8),
    @(41, 7, 8),
    @(42, 9, 9),
    @(43, 7, 7),
    @(44, 6, 6),
    @(45, 9, 9),
    @(46, 7, 7),
    @(47, 8, 8),
    @(48, 7, 7)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
